$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 5475
$ws.Range("I4").Value = 5475
$ws.Range("K4").Value = 5475
$ws.Range("M4").Value = -5361

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 100.53846
$ws.Range("I5").Value = 100.666664
$ws.Range("K5").Value = 100.666664
$ws.Range("M5").Value = 14.333336

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 195.85715
$ws.Range("I9").Value = 142.75
$ws.Range("K9").Value = 142.75
$ws.Range("M9").Value = 26.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 336.25
$ws.Range("I19").Value = 365
$ws.Range("K19").Value = 365
$ws.Range("M19").Value = -190

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 40295
$ws.Range("J68").Value = 40295
$ws.Range("L68").Value = 40295
$ws.Range("N68").Value = -41793

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 40295
$ws.Range("J71").Value = 40295
$ws.Range("L71").Value = 120885
$ws.Range("N71").Value = -128373

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 18145.5
$ws.Range("I141").Value = 860.6667
$ws.Range("J141").Value = 70000
$ws.Range("K141").Value = 2582.0001
$ws.Range("L141").Value = 210000
$ws.Range("M141").Value = 2597.9999
$ws.Range("N141").Value = -220360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1964.4286
$ws.Range("I2").Value = 2050.4
$ws.Range("J2").Value = 1749.5
$ws.Range("K2").Value = 2050.4
$ws.Range("L2").Value = 1749.5
$ws.Range("M2").Value = -1937.4
$ws.Range("N2").Value = -1975.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11932.9375
$ws.Range("I32").Value = 8268.727999999999
$ws.Range("K32").Value = 8268.727999999999
$ws.Range("M32").Value = -7981.727999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1964.4286
$ws.Range("I116").Value = 2050.4
$ws.Range("J116").Value = 1749.5
$ws.Range("K116").Value = 2050.4
$ws.Range("L116").Value = 1749.5
$ws.Range("M116").Value = 243.5999999999999
$ws.Range("N116").Value = -6337.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1964.4286
$ws.Range("I3").Value = 2050.4
$ws.Range("J3").Value = 1749.5
$ws.Range("K3").Value = 2050.4
$ws.Range("L3").Value = 1749.5
$ws.Range("M3").Value = -1936.4
$ws.Range("N3").Value = -1977.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3020.5715
$ws.Range("I86").Value = 2489
$ws.Range("K86").Value = 2489
$ws.Range("M86").Value = -1366

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3020.5715
$ws.Range("I89").Value = 2489
$ws.Range("K89").Value = 12445
$ws.Range("M89").Value = -6829

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 55000
$ws.Range("J80").Value = 55000
$ws.Range("L80").Value = 55000
$ws.Range("N80").Value = -57246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 55000
$ws.Range("J83").Value = 55000
$ws.Range("L83").Value = 165000
$ws.Range("N83").Value = -176232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2750
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -6246

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 2750
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -31232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I132").Value = 3046.8
$ws.Range("J132").Value = 2190
$ws.Range("K132").Value = 9140.400000000001
$ws.Range("L132").Value = 6570
$ws.Range("M132").Value = -6610.400000000001
$ws.Range("N132").Value = -11630

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 830
$ws.Range("I3").Value = 830
$ws.Range("K3").Value = 2490
$ws.Range("M3").Value = -2378

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 417.14285
$ws.Range("I4").Value = 153.33333
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 459.99999
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = -347.99999
$ws.Range("N4").Value = -6224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 57.555557
$ws.Range("I7").Value = 2.25
$ws.Range("K7").Value = 6.75
$ws.Range("M7").Value = 105.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 805.6
$ws.Range("I14").Value = 805.6
$ws.Range("K14").Value = 2416.8
$ws.Range("M14").Value = -2243.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2660.6155
$ws.Range("I131").Value = 1887.5555
$ws.Range("J131").Value = 4400
$ws.Range("K131").Value = 5662.666499999999
$ws.Range("L131").Value = 13200
$ws.Range("M131").Value = -622.6664999999994
$ws.Range("N131").Value = -23280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4225
$ws.Range("I14").Value = 4000
$ws.Range("J14").Value = 4300
$ws.Range("K14").Value = 4000
$ws.Range("L14").Value = 4300
$ws.Range("M14").Value = -3832
$ws.Range("N14").Value = -4636

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 9747.5
$ws.Range("J29").Value = 9747.5
$ws.Range("L29").Value = 9747.5
$ws.Range("N29").Value = -10327.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3332.6667
$ws.Range("I80").Value = 3249
$ws.Range("K80").Value = 3249
$ws.Range("M80").Value = -2251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3332.6667
$ws.Range("I83").Value = 3249
$ws.Range("K83").Value = 16245
$ws.Range("M83").Value = -11253

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 616.6667
$ws.Range("I113").Value = 616.6667
$ws.Range("K113").Value = 616.6667
$ws.Range("M113").Value = 1553.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 20004
$ws.Range("I3").Value = 20004
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 20004
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -19892
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 20004
$ws.Range("I15").Value = 20004
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 20004
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -19834
$ws.Range("N15").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7500
$ws.Range("I40").Value = 7000
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 7000
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -6864
$ws.Range("N40").Value = -8272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 49992
$ws.Range("J63").Value = 49992
$ws.Range("L63").Value = 49992
$ws.Range("N63").Value = -51490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 49992
$ws.Range("J66").Value = 49992
$ws.Range("L66").Value = 149976
$ws.Range("N66").Value = -157464

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2890.6365
$ws.Range("J68").Value = 3100
$ws.Range("L68").Value = 3100
$ws.Range("N68").Value = -4598

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2890.6365
$ws.Range("J71").Value = 3100
$ws.Range("L71").Value = 15500
$ws.Range("N71").Value = -22988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H102").Value = 75000
$ws.Range("J102").Value = 75000
$ws.Range("L102").Value = 75000
$ws.Range("N102").Value = -81490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 37745.668
$ws.Range("I132").Value = 35294.8
$ws.Range("K132").Value = 105884.4
$ws.Range("M132").Value = -103354.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 851.2
$ws.Range("I4").Value = 64
$ws.Range("K4").Value = 64
$ws.Range("M4").Value = 49

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 6669900
$ws.Range("I14").Value = 6669900
$ws.Range("K14").Value = 6669900
$ws.Range("M14").Value = -6669732

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3569.8
$ws.Range("I62").Value = 3462.25
$ws.Range("K62").Value = 3462.25
$ws.Range("M62").Value = -2838.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 3569.8
$ws.Range("I65").Value = 3462.25
$ws.Range("K65").Value = 17311.25
$ws.Range("M65").Value = -14191.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1000
$ws.Range("I96").Value = 800
$ws.Range("J96").Value = 1200
$ws.Range("K96").Value = 800
$ws.Range("L96").Value = 1200
$ws.Range("M96").Value = 573
$ws.Range("N96").Value = -3946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2311.125
$ws.Range("I136").Value = 2144.5386
$ws.Range("J136").Value = 3033
$ws.Range("K136").Value = 6433.6158
$ws.Range("L136").Value = 9099
$ws.Range("M136").Value = -3883.6158
$ws.Range("N136").Value = -14199
